$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text storage for numeric-looking strings without leaving
# a residual style index (set NumberFormat to Text, assign value, then restore the
# cell to the Normal style so no explicit s="n" attribute is written).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.728.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.418.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.58%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.35%  "

$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("E10").Value = "  +2.19%  "

$ws.Range("E11").Value = "  -2.08%  "

$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.845.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.653.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.425.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "329.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("E21").Value = "  -3.23%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.55%  "

$ws.Range("E24").Value = "  +3.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.46%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0774"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.83%  "

$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("E33").Value = "  +1.79%  "

$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  +4.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "314.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.409"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0968"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.30%  "

$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.579"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.70%  "

$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
